$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10 values
$ws.Range("C10").Value = "Corte Adulto"
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = "jose"

# Delete row 11 entirely (shift cells up)
$ws.Rows("11:11").Delete()
